# Updated output folder deploy:
#   - bump the IG "Date" / "Count" metadata on the Metadata sheet
#   - refresh the facility-type concepts on the Concepts sheet (6 concepts
#     now instead of 5: outreach-post, phc-center-l1/2/3, hospital, tertiary)

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Metadata sheet
# ----------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Date (plain text, no special handling needed)
$meta.Range("B8").Value = "2025-07-17T18:57:38+01:00"

# Count "5" -> "6" -- this column stores values as text, so force Text
# number format before writing the otherwise-numeric-looking string
$meta.Range("B23").NumberFormat = "@"
$meta.Range("B23").Value = "6"

# ----------------------------------------------------------------------
# Concepts sheet
# ----------------------------------------------------------------------
$concepts = $wb.Worksheets.Item("Concepts")

# Grow the table by one row (6 concepts -> 7 data rows incl. header).
# Copy the last existing row down first so the new row inherits the same
# formatting as the rest of the table.
$concepts.Range("A6:D6").Copy()
$concepts.Range("A7:D7").PasteSpecial(-4122)

# The "Level" column is always the text "1" -- keep it stored as text
# (rather than a number) for every data row, including the new one.
$concepts.Range("A2:A7").NumberFormat = "@"

# Row 2: clinic -> outreach-post
$concepts.Range("A2").Value = "1"
$concepts.Range("B2").Value = "outreach-post"
$concepts.Range("C2").Value = "Outreach Post"
$concepts.Range("D2").Value = "A temporary or mobile health service site established to deliver essential healthcare to underserved or remote populations."

# Row 3: hospital -> phc-center-l1
$concepts.Range("A3").Value = "1"
$concepts.Range("B3").Value = "phc-center-l1"
$concepts.Range("C3").Value = "PHC Center Level 1"
$concepts.Range("D3").Value = "A basic primary healthcare facility providing preventive and promotive services, often staffed by community health workers."

# Row 4: health-post -> phc-center-l2
$concepts.Range("A4").Value = "1"
$concepts.Range("B4").Value = "phc-center-l2"
$concepts.Range("C4").Value = "PHC Center Level 2"
$concepts.Range("D4").Value = "A primary healthcare facility with limited diagnostic and treatment services, typically staffed by nurses, CHEWs, or junior medical personnel."

# Row 5: chc -> phc-center-l3
$concepts.Range("A5").Value = "1"
$concepts.Range("B5").Value = "phc-center-l3"
$concepts.Range("C5").Value = "PHC Center Level 3"
$concepts.Range("D5").Value = "An advanced primary healthcare facility offering comprehensive outpatient care, minor/light procedures, and referral support, often with a resident medical officer."

# Row 6: tertiary -> hospital (now "Secondary Hospital")
$concepts.Range("A6").Value = "1"
$concepts.Range("B6").Value = "hospital"
$concepts.Range("C6").Value = "Secondary Hospital"
$concepts.Range("D6").Value = "A referral facility providing specialized care, inpatient services, and emergency response, staffed by medical officers and specialists."

# Row 7 (new): tertiary (now "Tertiary Facility")
$concepts.Range("A7").Value = "1"
$concepts.Range("B7").Value = "tertiary"
$concepts.Range("C7").Value = "Tertiary Facility"
$concepts.Range("D7").Value = "A highly specialized hospital offering advanced diagnostic, therapeutic, and surgical care, often affiliated with teaching or research institutions."

$concepts.Range("A1").Select()
